$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column (H1), reusing the same formatting as the
# neighboring header cells (e.g. G1: bold, bordered, centered).
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add the corresponding numeric value for row 2 (H2)
$ws.Range("H2").Value = 1
